$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet: Last Updated timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 01:11 AM"

# --- Update "Top Gainers" sheet: shift rows 36-75 up (drop old row 36 / SKMEGGPROD)
#     and append a brand-new row 76 (CGPOWER) ---
$gainers = $wb.Worksheets.Item("Top Gainers")

$rows = @(
    @("PROZONER",   4.9921,  15.7468,  36.095),
    @("STALLION",   4.9914,  -5.2229,  21.4391),
    @("SURYAROSNI", 4.9831,  11.386,   3.0213),
    @("DATAMATICS", 4.9005,  7.3152,   15.7298),
    @("UTKARSHBNK", 4.8768,  -5.8959,  -2.6215),
    @("FILATEX",    4.8689,  10.274,   26.0027),
    @("INDOTHAI",   4.8064,  4.5349,   43.748),
    @("SANDUMA",    4.593,   2.1405,   30.2813),
    @("LLOYDSENT",  4.5646,  1.8339,   11.234),
    @("STAR",       4.5025,  4.4319,   3.662),
    @("RECLTD",     4.4992,  3.4756,   3.4062),
    @("NBCC",       4.4511,  3.1605,   7.6018),
    @("GPPL",       4.4154,  3.4073,   5.0497),
    @("BIL",        4.3654,  9.122199999999999, -0.3203),
    @("HUDCO",      4.3201,  3.8924,   5.3884),
    @("SGMART",     4.2736,  8.258900000000001, 2.5381),
    @("MRPL",       4.2642,  9.7103,   20.0542),
    @("JKIL",       4.1372,  2.9463,   1.7584),
    @("SAMBHV",     4.1349,  2.624,    5.167),
    @("SAPPHIRE",   4.1265,  1.7633,   -0.7999000000000001),
    @("PVRINOX",    4.1118,  6.2102,   14.707),
    @("KERNEX",     3.9981,  7.4592,   27.1054),
    @("SUNFLAG",    3.997,   4.333,    4.6312),
    @("CMSINFO",    3.9096,  2.6872,   2.8935),
    @("GMBREW",     3.8999,  -0.53,    79.029),
    @("APARINDS",   3.8924,  8.3414,   15.5876),
    @("HITECHGEAR", 3.8587,  1.1486,   9.9254),
    @("NPST",       3.8509,  -2.0059,  -3.5057),
    @("ORIENTTECH", 3.827,   0.5247000000000001, 32.6784),
    @("ICRA",       3.7985,  4.4793,   2.8828),
    @("SALASAR",    3.7935,  4.7872,   11.0485),
    @("DCW",        3.7544,  2.3219,   -3.9753),
    @("RHETAN",     3.754,   4.178,    6.549),
    @("HINDPETRO",  3.6935,  6.9335,   5.7397),
    @("BHARTIHEXA", 3.6718,  7.0877,   15.3332),
    @("HLEGLAS",    3.659,   8.115500000000001, 27.1239),
    @("RHIM",       3.6544,  3.2276,   5.1826),
    @("SHK",        3.6347,  2.388,    -1.932),
    @("BCLIND",     3.6271,  2.2945,   0.1728),
    @("MUKANDLTD",  3.6133,  11.9685,  9.550800000000001),
    @("CGPOWER",    3.6125,  3.4192,   1.0325)
)

$startRow = 36
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]
    $gainers.Cells.Item($r, 2).Value = $item[0]
    $gainers.Cells.Item($r, 3).Value = $item[1]
    $gainers.Cells.Item($r, 4).Value = $item[2]
    $gainers.Cells.Item($r, 5).Value = $item[3]
}
